$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 394.23334
$ws.Range("I92").Value = 276.57144
$ws.Range("K92").Value = 276.57144
$ws.Range("M92").Value = 971.4285600000001
$ws.Range("H96").Value = 1487.3077
$ws.Range("I96").Value = 1820.75
$ws.Range("J96").Value = 953.8
$ws.Range("K96").Value = 5462.25
$ws.Range("L96").Value = 2861.4
$ws.Range("M96").Value = -4089.25
$ws.Range("N96").Value = -5607.4
$ws.Range("H103").Value = 347820.12
$ws.Range("I103").Value = 926512
$ws.Range("J103").Value = 605
$ws.Range("K103").Value = 2779536
$ws.Range("L103").Value = 1815
$ws.Range("M103").Value = -2778950
$ws.Range("N103").Value = -2987
$ws.Range("H138").Value = 3007
$ws.Range("I138").Value = 1020.8205
$ws.Range("J138").Value = 4655.1064
$ws.Range("K138").Value = 3062.4615
$ws.Range("L138").Value = 13965.3192
$ws.Range("M138").Value = 2077.5385
$ws.Range("N138").Value = -24245.3192
$ws.Range("H141").Value = 1953.1428
$ws.Range("I141").Value = 1264.35
$ws.Range("J141").Value = 5014.4443
$ws.Range("K141").Value = 3793.05
$ws.Range("L141").Value = 15043.3329
$ws.Range("M141").Value = 1386.95
$ws.Range("N141").Value = -25403.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2351.6711
$ws.Range("I63").Value = 2303.7246
$ws.Range("J63").Value = 2824.2856
$ws.Range("K63").Value = 2303.7246
$ws.Range("L63").Value = 2824.2856
$ws.Range("M63").Value = -1617.7246
$ws.Range("N63").Value = -4196.2856
$ws.Range("H66").Value = 2351.6711
$ws.Range("I66").Value = 2303.7246
$ws.Range("J66").Value = 2824.2856
$ws.Range("K66").Value = 11518.623
$ws.Range("L66").Value = 14121.428
$ws.Range("M66").Value = -8086.623
$ws.Range("N66").Value = -20985.428
$ws.Range("H109").Value = 59800
$ws.Range("J109").Value = 59800
$ws.Range("L109").Value = 59800
$ws.Range("N109").Value = -62574
$ws.Range("H122").Value = 1902.1538
$ws.Range("I122").Value = 1148.24
$ws.Range("J122").Value = 3248.4285
$ws.Range("K122").Value = 3444.72
$ws.Range("L122").Value = 9745.2855
$ws.Range("M122").Value = -994.7200000000003
$ws.Range("N122").Value = -14645.2855
$ws.Range("H132").Value = 1917.3793
$ws.Range("I132").Value = 1698
$ws.Range("J132").Value = 2276.3635
$ws.Range("K132").Value = 5094
$ws.Range("L132").Value = 6829.0905
$ws.Range("M132").Value = -2564
$ws.Range("N132").Value = -11889.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 20684
$ws.Range("J108").Value = 20684
$ws.Range("L108").Value = 20684
$ws.Range("N108").Value = -28364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3646.182
$ws.Range("I62").Value = 3567
$ws.Range("J62").Value = 3741.2
$ws.Range("K62").Value = 3567
$ws.Range("L62").Value = 3741.2
$ws.Range("M62").Value = -2943
$ws.Range("N62").Value = -4989.2
$ws.Range("H65").Value = 3646.182
$ws.Range("I65").Value = 3567
$ws.Range("J65").Value = 3741.2
$ws.Range("K65").Value = 17835
$ws.Range("L65").Value = 18706
$ws.Range("M65").Value = -14715
$ws.Range("N65").Value = -24946
$ws.Range("H132").Value = 1847.625
$ws.Range("I132").Value = 1462.25
$ws.Range("J132").Value = 2618.375
$ws.Range("K132").Value = 4386.75
$ws.Range("L132").Value = 7855.125
$ws.Range("M132").Value = -1856.75
$ws.Range("N132").Value = -12915.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 619.2381
$ws.Range("I5").Value = 548.0909
$ws.Range("J5").Value = 697.5
$ws.Range("K5").Value = 1644.2727
$ws.Range("L5").Value = 2092.5
$ws.Range("M5").Value = -1532.2727
$ws.Range("N5").Value = -2316.5
$ws.Range("H69").Value = 1163.3636
$ws.Range("I69").Value = 966.1667
$ws.Range("J69").Value = 1400
$ws.Range("K69").Value = 2898.5001
$ws.Range("L69").Value = 4200
$ws.Range("M69").Value = -2087.5001
$ws.Range("N69").Value = -5822
$ws.Range("H72").Value = 1163.3636
$ws.Range("I72").Value = 966.1667
$ws.Range("J72").Value = 1400
$ws.Range("K72").Value = 8695.5003
$ws.Range("L72").Value = 12600
$ws.Range("M72").Value = -4639.5003
$ws.Range("N72").Value = -20712
$ws.Range("H98").Value = 183.83333
$ws.Range("I98").Value = 160.6
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 481.8
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = 1016.2
$ws.Range("N98").Value = -3896
$ws.Range("H122").Value = 4167420.8
$ws.Range("I122").Value = 605.5
$ws.Range("J122").Value = 12501051
$ws.Range("K122").Value = 5449.5
$ws.Range("L122").Value = 112509459
$ws.Range("M122").Value = -2999.5
$ws.Range("N122").Value = -112514359
$ws.Range("H135").Value = 619.2381
$ws.Range("I135").Value = 548.0909
$ws.Range("J135").Value = 697.5
$ws.Range("K135").Value = 4932.8181
$ws.Range("L135").Value = 6277.5
$ws.Range("M135").Value = -2397.8181
$ws.Range("N135").Value = -11347.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 13253.667
$ws.Range("J39").Value = 13253.667
$ws.Range("L39").Value = 13253.667
$ws.Range("N39").Value = -14317.667
$ws.Range("H111").Value = 18000
$ws.Range("J111").Value = 18000
$ws.Range("L111").Value = 18000
$ws.Range("N111").Value = -24134
$ws.Range("H132").Value = 6830.6665
$ws.Range("I132").Value = 7934.6665
$ws.Range("J132").Value = 3518.6667
$ws.Range("K132").Value = 23803.9995
$ws.Range("L132").Value = 10556.0001
$ws.Range("M132").Value = -21273.9995
$ws.Range("N132").Value = -15616.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 9999.5
$ws.Range("I4").Value = 9999
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 9999
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -9886
$ws.Range("N4").Value = -10226
$ws.Range("H28").Value = 9999.5
$ws.Range("I28").Value = 9999
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 9999
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = -9767
$ws.Range("N28").Value = -10464
$ws.Range("H37").Value = 9999.5
$ws.Range("I37").Value = 9999
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 9999
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -9892
$ws.Range("N37").Value = -10214
$ws.Range("H46").Value = 1353.2222
$ws.Range("I46").Value = 1199.5
$ws.Range("J46").Value = 1660.6666
$ws.Range("K46").Value = 1199.5
$ws.Range("L46").Value = 1660.6666
$ws.Range("M46").Value = -1011.5
$ws.Range("N46").Value = -2036.6666
$ws.Range("H55").Value = 224.58333
$ws.Range("I55").Value = 180
$ws.Range("J55").Value = 298.8889
$ws.Range("K55").Value = 180
$ws.Range("L55").Value = 298.8889
$ws.Range("M55").Value = -7
$ws.Range("N55").Value = -644.8888999999999
$ws.Range("H98").Value = 44800
$ws.Range("J98").Value = 44800
$ws.Range("L98").Value = 44800
$ws.Range("N98").Value = -50790
$ws.Range("H132").Value = 2946.25
$ws.Range("I132").Value = 1963.4117
$ws.Range("J132").Value = 3825.6316
$ws.Range("K132").Value = 5890.2351
$ws.Range("L132").Value = 11476.8948
$ws.Range("M132").Value = -3360.2351
$ws.Range("N132").Value = -16536.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 914136.4
$ws.Range("J2").Value = 914136.4
$ws.Range("L2").Value = 914136.4
$ws.Range("N2").Value = -914360.4
$ws.Range("H62").Value = 3537.2307
$ws.Range("I62").Value = 3136.2
$ws.Range("J62").Value = 3787.875
$ws.Range("K62").Value = 3136.2
$ws.Range("L62").Value = 3787.875
$ws.Range("M62").Value = -2512.2
$ws.Range("N62").Value = -5035.875
$ws.Range("H65").Value = 3537.2307
$ws.Range("I65").Value = 3136.2
$ws.Range("J65").Value = 3787.875
$ws.Range("K65").Value = 15681
$ws.Range("L65").Value = 18939.375
$ws.Range("M65").Value = -12561
$ws.Range("N65").Value = -25179.375
$ws.Range("H132").Value = 1431.2858
$ws.Range("I132").Value = 1153.4117
$ws.Range("J132").Value = 1860.7273
$ws.Range("K132").Value = 3460.2351
$ws.Range("L132").Value = 5582.1819
$ws.Range("M132").Value = -930.2351000000003
$ws.Range("N132").Value = -10642.1819
